# Reproduces the upload diff for "Дагтрио.xlsx":
#  - Row 53 on the "Особые события" sheet gets a bottom-border style
#    (it becomes the last row of the "us3102.ssb" group).
#  - Six new rows (54-59) are appended for a new "us2302.ssb" group,
#    each with an English line (C), Russian translation (D) and the
#    "converted" cipher text (E); row 54 also carries the script file
#    name in column A. These new rows reuse the plain (non-bordered)
#    style used by the rows above them.
#  - The selected cell moves to D58 to reflect the new last-edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Особые события")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Apply the plain (non-bordered) formatting to the new rows 54-59
#    first, copying it from the still-unmodified rows 49-52 (style
#    s="4"/"5"). Row 53 is restyled to the bordered variant later, so
#    it must not be used as a format source beforehand.
# ---------------------------------------------------------------------
$ws.Range("A49:E49").Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)

$ws.Range("A49:E49").Copy()
$ws.Range("A55:E55").PasteSpecial(-4122)

$ws.Range("A49:E49").Copy()
$ws.Range("A56:E56").PasteSpecial(-4122)

$ws.Range("A49:E49").Copy()
$ws.Range("A57:E57").PasteSpecial(-4122)

$ws.Range("A49:E49").Copy()
$ws.Range("A58:E58").PasteSpecial(-4122)

$ws.Range("A49:E49").Copy()
$ws.Range("A59:E59").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Fill in the values. The order below (all of column C, then A54 +
#    column D, then column E) matches the order new strings were
#    appended to the shared-string table in the original edit.
# ---------------------------------------------------------------------
$ws.Range("C54").Value = " Gwah![K] Gwargh![K] Gwagagah!"
$ws.Range("C55").Value = " We thought we were a goner!"
$ws.Range("C56").Value = " While training here...[K]the waves\ncarried us off!"
$ws.Range("C57").Value = " For months, we drifted on the\nwaves, cut off from all contact…"
$ws.Range("C58").Value = " ...In all truth, that is really an\nexaggeration.[K] We were merely splashed by\nwaves on the shore…"
$ws.Range("C59").Value = " While we floundered...[K]we were\nthankfully rescued by these kindly strangers."

$ws.Range("A54").Value = "SCRIPT/D01P11A/us2302.ssb"
$ws.Range("D54").Value = " Гвах![K] Гвааарх![K] Гвагахах!"
$ws.Range("D55").Value = " Мы думали, нам конец!"
$ws.Range("D56").Value = " Мы тренировались здесь...[K]\nИ волны унесли нас!"
$ws.Range("D57").Value = " Месяцами мы дрейфовали по волнам,\nотрезанными от всех..."
$ws.Range("D58").Value = " ...По правде говоря, это несколько\nпреувеличено.[K] Нас лишь немножечко\nнамочило пляжной волной..."
$ws.Range("D59").Value = " Пока мы барахтались...[K] Нас спасли\nэти добрые незнакомцы."

$ws.Range("E54").Value = " Ãâàö![K] Ãâàààñö![K] Ãâàãàöàö!"
$ws.Range("E55").Value = " Íú äôíàìé, îàí ëïîåø!"
$ws.Range("E56").Value = " Íú óñåîéñïâàìéòû èäåòû...[K]\nÉ âïìîú ôîåòìé îàò!"
$ws.Range("E57").Value = " Íåòÿøàíé íú äñåêõïâàìé ðï âïìîàí,\nïóñåèàîîúíé ïó âòåö..."
$ws.Range("E58").Value = " ...Ðï ðñàâäå ãïâïñÿ, üóï îåòëïìûëï\nðñåôâåìéœåîï.[K] Îàò ìéšû îåíîïçåœëï\nîàíïœéìï ðìÿçîïê âïìîïê..."
$ws.Range("E59").Value = " Ðïëà íú áàñàöóàìéòû...[K] Îàò òðàòìé\nüóé äïáñúå îåèîàëïíøú."

# ---------------------------------------------------------------------
# 3) Row-specific numbers (column B).
# ---------------------------------------------------------------------
$ws.Range("B54").Value = 26
$ws.Range("B55").Value = 36
$ws.Range("B56").Value = 39
$ws.Range("B57").Value = 42
$ws.Range("B58").Value = 46
$ws.Range("B59").Value = 54

# ---------------------------------------------------------------------
# 4) Re-style row 53 to the bordered "group end" variant, copying the
#    formatting from row 48 (the previous group's closing row).
# ---------------------------------------------------------------------
$ws.Range("A48:E48").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5) Update the selection to match the new last-edited cell.
# ---------------------------------------------------------------------
$ws.Range("D58").Select()
